$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.708798885345459
$ws.Range("E2").Value = 1448.05328566901
$ws.Range("F2").Value = 0.06734148368525826
$ws.Range("G2").Value = 0.05505040562410878
$ws.Range("H2").Value = 0.04792491229542371
$ws.Range("I2").Value = 0.04236825748483357
$ws.Range("J2").Value = 0.03807524420185895
$ws.Range("K2").Value = 0.03454229328210902
$ws.Range("L2").Value = 0.03213501245672883
$ws.Range("M2").Value = 0.03184823135405263
$ws.Range("N2").Value = 0.03092115310218088
$ws.Range("O2").Value = 0.03059610685946856
$ws.Range("P2").Value = 0.02988784642409726
$ws.Range("Q2").Value = 0.02942597165874971
$ws.Range("R2").Value = 0.02915170243356777
$ws.Range("S2").Value = 0.02885458996040594
$ws.Range("T2").Value = 0.02878267415235739
$ws.Range("U2").Value = 0.02854347143364555
$ws.Range("V2").Value = 0.02839973429147171
$ws.Range("W2").Value = 0.02839973429147171
$ws.Range("X2").Value = 0.02831114312969657
$ws.Range("Y2").Value = 0.02822715956469804
$ws.Range("C3").Value = 0.7948737144470215
$ws.Range("E3").Value = 1489.904564869592
$ws.Range("F3").Value = 0.06810240487318925
$ws.Range("G3").Value = 0.05517304379635206
$ws.Range("H3").Value = 0.0449984857945104
$ws.Range("I3").Value = 0.04145986882199669
$ws.Range("J3").Value = 0.03902636456940533
$ws.Range("K3").Value = 0.03723394995561809
$ws.Range("L3").Value = 0.03570663423462161
$ws.Range("M3").Value = 0.03350323711037784
$ws.Range("N3").Value = 0.0327464322978905
$ws.Range("O3").Value = 0.03209542636751661
$ws.Range("P3").Value = 0.0310561342390011
$ws.Range("Q3").Value = 0.02967495056841576
$ws.Range("R3").Value = 0.02967495056841576
$ws.Range("S3").Value = 0.02945533037196667
$ws.Range("T3").Value = 0.02940175733040556
$ws.Range("U3").Value = 0.02912829677037918
$ws.Range("V3").Value = 0.02912829677037918
$ws.Range("W3").Value = 0.02912829677037918
$ws.Range("X3").Value = 0.02906434096151064
$ws.Range("Y3").Value = 0.02904297397406611
$ws.Range("C4").Value = 0.7354607582092285
$ws.Range("E4").Value = 1541.152179596231
$ws.Range("F4").Value = 0.06775691073841195
$ws.Range("G4").Value = 0.05434322035179644
$ws.Range("H4").Value = 0.04814376526188227
$ws.Range("I4").Value = 0.04525903684008887
$ws.Range("J4").Value = 0.04038739542168356
$ws.Range("K4").Value = 0.0370802751993872
$ws.Range("L4").Value = 0.03607167413999854
$ws.Range("M4").Value = 0.03448418730544573
$ws.Range("N4").Value = 0.03343284984446292
$ws.Range("O4").Value = 0.03247872715987322
$ws.Range("P4").Value = 0.03192364383010562
$ws.Range("Q4").Value = 0.03158190835925204
$ws.Range("R4").Value = 0.03138525698268505
$ws.Range("S4").Value = 0.03093965111721659
$ws.Range("T4").Value = 0.03060720307819696
$ws.Range("U4").Value = 0.03045851068969805
$ws.Range("V4").Value = 0.0302786638459095
$ws.Range("W4").Value = 0.03024104488340618
$ws.Range("X4").Value = 0.03011668876375502
$ws.Range("Y4").Value = 0.03004195281863997
$ws.Range("C5").Value = 0.7031199932098389
$ws.Range("E5").Value = 1497.186824226579
$ws.Range("F5").Value = 0.0672456839988628
$ws.Range("G5").Value = 0.05657568496370059
$ws.Range("H5").Value = 0.0477458644051469
$ws.Range("I5").Value = 0.04228088860076793
$ws.Range("J5").Value = 0.03933476695716807
$ws.Range("K5").Value = 0.03692756930129316
$ws.Range("L5").Value = 0.03454001101886991
$ws.Range("M5").Value = 0.03347246870684517
$ws.Range("N5").Value = 0.03276766262713533
$ws.Range("O5").Value = 0.0316495150884161
$ws.Range("P5").Value = 0.03083621275591808
$ws.Range("Q5").Value = 0.03034223198578209
$ws.Range("R5").Value = 0.0299177808342243
$ws.Range("S5").Value = 0.02971942387621307
$ws.Range("T5").Value = 0.02971836810056327
$ws.Range("U5").Value = 0.0296236739641736
$ws.Range("V5").Value = 0.02952765577235657
$ws.Range("W5").Value = 0.02925784386087629
$ws.Range("X5").Value = 0.02925784386087629
$ws.Range("Y5").Value = 0.02918492834749666
$ws.Range("C6").Value = 0.7187392711639404
$ws.Range("E6").Value = 1545.066929621027
$ws.Range("F6").Value = 0.06762117489747681
$ws.Range("G6").Value = 0.05345745275471724
$ws.Range("H6").Value = 0.0472760459665636
$ws.Range("I6").Value = 0.04302408632300653
$ws.Range("J6").Value = 0.03997455093440801
$ws.Range("K6").Value = 0.03673463312149364
$ws.Range("L6").Value = 0.0344684096322875
$ws.Range("M6").Value = 0.03331944482693355
$ws.Range("N6").Value = 0.0329920573656258
$ws.Range("O6").Value = 0.03236778302992665
$ws.Range("P6").Value = 0.03232901676534485
$ws.Range("Q6").Value = 0.03163937897794613
$ws.Range("R6").Value = 0.0312920460443203
$ws.Range("S6").Value = 0.0310526447490998
$ws.Range("T6").Value = 0.03072483742442202
$ws.Range("U6").Value = 0.03056705857797291
$ws.Range("V6").Value = 0.03041184994264001
$ws.Range("W6").Value = 0.03032951668868616
$ws.Range("X6").Value = 0.03017097100123558
$ws.Range("Y6").Value = 0.03011826373530267
$ws.Range("C7").Value = 0.7187497615814209
$ws.Range("E7").Value = 1533.58244787914
$ws.Range("F7").Value = 0.06731517371662961
$ws.Range("G7").Value = 0.05498301897314078
$ws.Range("H7").Value = 0.04782023718856787
$ws.Range("I7").Value = 0.04514663880210839
$ws.Range("J7").Value = 0.04131143266402484
$ws.Range("K7").Value = 0.03836842141566298
$ws.Range("L7").Value = 0.03654387252830766
$ws.Range("M7").Value = 0.0354126450905902
$ws.Range("N7").Value = 0.03411159162300496
$ws.Range("O7").Value = 0.033659050097132
$ws.Range("P7").Value = 0.03297813068918972
$ws.Range("Q7").Value = 0.03239670621568314
$ws.Range("R7").Value = 0.03117620842964664
$ws.Range("S7").Value = 0.03072741331048395
$ws.Range("T7").Value = 0.03045261616871347
$ws.Range("U7").Value = 0.03043917275957136
$ws.Range("V7").Value = 0.03017431874748723
$ws.Range("W7").Value = 0.03008787394197295
$ws.Range("X7").Value = 0.02990675865036339
$ws.Range("Y7").Value = 0.0298943946954998
$ws.Range("C8").Value = 0.7031259536743164
$ws.Range("E8").Value = 1523.053553910388
$ws.Range("F8").Value = 0.06501018666108467
$ws.Range("G8").Value = 0.0545228238009107
$ws.Range("H8").Value = 0.0479539374492574
$ws.Range("I8").Value = 0.04580723083442008
$ws.Range("J8").Value = 0.04268867611053955
$ws.Range("K8").Value = 0.03946883614934294
$ws.Range("L8").Value = 0.03672258647081109
$ws.Range("M8").Value = 0.03476627921585379
$ws.Range("N8").Value = 0.03419691497416947
$ws.Range("O8").Value = 0.0326067168873132
$ws.Range("P8").Value = 0.03191753174263255
$ws.Range("Q8").Value = 0.03133585307578327
$ws.Range("R8").Value = 0.03101515850373525
$ws.Range("S8").Value = 0.03035872719602457
$ws.Range("T8").Value = 0.03035872719602457
$ws.Range("U8").Value = 0.03021754591801571
$ws.Range("V8").Value = 0.03005012499922026
$ws.Range("W8").Value = 0.02990032753930467
$ws.Range("X8").Value = 0.02976148800669126
$ws.Range("Y8").Value = 0.02968915309766837
$ws.Range("C9").Value = 0.7031211853027344
$ws.Range("E9").Value = 1514.261986118221
$ws.Range("F9").Value = 0.06676943473561732
$ws.Range("G9").Value = 0.05264553125347945
$ws.Range("H9").Value = 0.04622361771963186
$ws.Range("I9").Value = 0.04406999175272695
$ws.Range("J9").Value = 0.03843668393532098
$ws.Range("K9").Value = 0.03672302073926133
$ws.Range("L9").Value = 0.03495547076223143
$ws.Range("M9").Value = 0.03408338568551982
$ws.Range("N9").Value = 0.03329402667793729
$ws.Range("O9").Value = 0.03262107783706859
$ws.Range("P9").Value = 0.03143566260652664
$ws.Range("Q9").Value = 0.03143566260652664
$ws.Range("R9").Value = 0.03122779820159909
$ws.Range("S9").Value = 0.03066628741045385
$ws.Range("T9").Value = 0.03057366669322176
$ws.Range("U9").Value = 0.03010439244833343
$ws.Range("V9").Value = 0.02980353292195573
$ws.Range("W9").Value = 0.02980353292195573
$ws.Range("X9").Value = 0.02961753363838246
$ws.Range("Y9").Value = 0.02951777750717779
$ws.Range("C10").Value = 0.7187502384185791
$ws.Range("E10").Value = 1505.778464952426
$ws.Range("F10").Value = 0.0678416358667131
$ws.Range("G10").Value = 0.05645771709134453
$ws.Range("H10").Value = 0.04738540092135047
$ws.Range("I10").Value = 0.04043907830815214
$ws.Range("J10").Value = 0.03854609583932683
$ws.Range("K10").Value = 0.03685117388999191
$ws.Range("L10").Value = 0.0348694310253166
$ws.Range("M10").Value = 0.03363071655283082
$ws.Range("N10").Value = 0.03313617387564987
$ws.Range("O10").Value = 0.03233395056417573
$ws.Range("P10").Value = 0.03156898208031031
$ws.Range("Q10").Value = 0.03118527488470511
$ws.Range("R10").Value = 0.03026709342830818
$ws.Range("S10").Value = 0.03020777005270174
$ws.Range("T10").Value = 0.02982793391154151
$ws.Range("U10").Value = 0.02982793391154151
$ws.Range("V10").Value = 0.02960267994422603
$ws.Range("W10").Value = 0.02956679978102322
$ws.Range("X10").Value = 0.02946298024877584
$ws.Range("Y10").Value = 0.02935240672421881
$ws.Range("C11").Value = 0.7031300067901611
$ws.Range("E11").Value = 1533.507702813768
$ws.Range("F11").Value = 0.06787085321713363
$ws.Range("G11").Value = 0.05607279188414514
$ws.Range("H11").Value = 0.04708460183074994
$ws.Range("I11").Value = 0.04456994579466091
$ws.Range("J11").Value = 0.04172841837812251
$ws.Range("K11").Value = 0.03826035237248345
$ws.Range("L11").Value = 0.03671104606460574
$ws.Range("M11").Value = 0.03576829159578882
$ws.Range("N11").Value = 0.03439273434707859
$ws.Range("O11").Value = 0.03278351632076679
$ws.Range("P11").Value = 0.03214940820282364
$ws.Range("Q11").Value = 0.03196132506912719
$ws.Range("R11").Value = 0.03161442803019304
$ws.Range("S11").Value = 0.03105169659336472
$ws.Range("T11").Value = 0.03068357354758817
$ws.Range("U11").Value = 0.03025400775556457
$ws.Range("V11").Value = 0.03023523669255527
$ws.Range("W11").Value = 0.03010335496342829
$ws.Range("X11").Value = 0.02996353699889102
$ws.Range("Y11").Value = 0.02989293767668162
